# Freelancer weekday dayshift changed to 0930-1830
# Updates every weekday (Mon-Fri) "10-19" shift cell to "0930-1830".
# Weekend "10-19" cells are intentionally left untouched.
# Also corrects an anomalous Daisy/Tak entry on 01/04/2025 (row 17)
# from "7-16" back to the standard "13-22".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToUpdate = @(
    "G2", "B3", "G4", "D5", "E6",
    "C9", "C10", "G11", "F12", "G13",
    "B16", "G17", "G18", "D19", "C20",
    "D23", "B24", "D25", "D26", "B27",
    "D30", "D31"
)

foreach ($addr in $cellsToUpdate) {
    $ws.Range($addr).Value = "0930-1830"
}

$ws.Range("H17").Value = "13-22"
$ws.Range("I17").Value = "13-22"
